# Edit script for codigos/comparacao_resultados.xlsx
# Renames the "docente"/"preferencia" group-prefixed headers to prefix-first
# naming (docente_dep -> dep_docente, etc.) and refreshes the matched
# teacher/preference assignment data in columns C:J (re-run of the
# comparison script against updated inputs), then reapplies the
# auto-fit column widths and restores the working selection/scroll
# position, mirroring the resave performed in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (row 1), columns C:J ---
$ws.Range("C1").Value = "dep_docente"
$ws.Range("D1").Value = "dep_preferencia"
$ws.Range("E1").Value = "pli_docente"
$ws.Range("F1").Value = "pli_preferencia"
$ws.Range("G1").Value = "aco_docente"
$ws.Range("H1").Value = "aco_preferencia"
$ws.Range("I1").Value = "ag_docente"
$ws.Range("J1").Value = "ag_preferencia"

# --- Refresh aco_docente/aco_preferencia (G:H) and ag_docente/ag_preferencia (I:J) values ---
$ws.Range("G2").Value = "DONALD MATTHEW PIANTO"
$ws.Range("I2").Value = "FELIPE SOUSA QUINTINO"
$ws.Range("J2").Value = 2
$ws.Range("I3").Value = "JAMES MATOS SAMPAIO"
$ws.Range("J3").Value = 3
$ws.Range("G6").Value = "TEREZINHA KESSIA DE ASSIS RIBEIRO"
$ws.Range("I7").Value = "EDUARDO YOSHIO NAKANO"
$ws.Range("G8").Value = "LUIS GUSTAVO DO AMARAL VINHA"
$ws.Range("I8").Value = "CIRA ETHEOWALDA GUEVARA OTINIANO"
$ws.Range("J8").Value = 1
$ws.Range("G9").Value = "JOSE ANGELO BELLONI"
$ws.Range("I9").Value = "DEMERSON ANDRE POLLI"
$ws.Range("G10").Value = "ERITON BARROS DOS SANTOS"
$ws.Range("I10").Value = "ANA MARIA NOGALES VASCONCELOS"
$ws.Range("I11").Value = "LUIS GUSTAVO DO AMARAL VINHA"
$ws.Range("G12").Value = "ANA MARIA NOGALES VASCONCELOS"
$ws.Range("I12").Value = "NICOLLAS STEFAN SOARES DA COSTA"
$ws.Range("I13").Value = "DONALD MATTHEW PIANTO"
$ws.Range("I14").Value = "NICOLLAS STEFAN SOARES DA COSTA"
$ws.Range("G15").Value = "JOSE ANGELO BELLONI"
$ws.Range("I15").Value = "GUSTAVO LEONEL GILARDONI AVALLE"
$ws.Range("G16").Value = "ANA MARIA NOGALES VASCONCELOS"
$ws.Range("I16").Value = "LUIS GUSTAVO DO AMARAL VINHA"
$ws.Range("G17").Value = "NICOLLAS STEFAN SOARES DA COSTA"
$ws.Range("I17").Value = "ERITON BARROS DOS SANTOS"
$ws.Range("I18").Value = "EDUARDO YOSHIO NAKANO"
$ws.Range("G19").Value = "JOANLISE MARCO DE LEON ANDRADE"
$ws.Range("I19").Value = "PETER ZORNIG"
$ws.Range("G20").Value = "ROBERTO VILA GABRIEL"
$ws.Range("I20").Value = "JAMES MATOS SAMPAIO"
$ws.Range("G21").Value = "JOANLISE MARCO DE LEON ANDRADE"
$ws.Range("I21").Value = "ROBERTO VILA GABRIEL"
$ws.Range("G22").Value = "JOSE AUGUSTO FIORUCCI"
$ws.Range("I22").Value = "GUSTAVO LEONEL GILARDONI AVALLE"
$ws.Range("G23").Value = "GUSTAVO LEONEL GILARDONI AVALLE"
$ws.Range("I23").Value = "JOSE ANGELO BELLONI"
$ws.Range("G24").Value = "RAUL YUKIHIRO MATSUSHITA"
$ws.Range("I24").Value = "ERITON BARROS DOS SANTOS"
$ws.Range("G25").Value = "PETER ZORNIG"
$ws.Range("I25").Value = "ANTONIO EDUARDO GOMES"
$ws.Range("G26").Value = "GUSTAVO LEONEL GILARDONI AVALLE"
$ws.Range("I26").Value = "PETER ZORNIG"
$ws.Range("G27").Value = "ANDRE LUIZ FERNANDES CANCADO"
$ws.Range("G28").Value = "DONALD MATTHEW PIANTO"
$ws.Range("I28").Value = "EDUARDO YOSHIO NAKANO"
$ws.Range("G29").Value = "RAUL YUKIHIRO MATSUSHITA"
$ws.Range("I29").Value = "DEMERSON ANDRE POLLI"
$ws.Range("J29").Value = 3
$ws.Range("G30").Value = "GUSTAVO LEONEL GILARDONI AVALLE"
$ws.Range("I30").Value = "ALAN RICARDO DA SILVA"
$ws.Range("J30").Value = 3
$ws.Range("I31").Value = "MARIA TERESA LEAO COSTA"
$ws.Range("G32").Value = "ANTONIO EDUARDO GOMES"
$ws.Range("I32").Value = "ANTONIO EDUARDO GOMES"
$ws.Range("G37").Value = "ROBERTO VILA GABRIEL"
$ws.Range("G38").Value = "JAMES MATOS SAMPAIO"
$ws.Range("G40").Value = "CIRA ETHEOWALDA GUEVARA OTINIANO"
$ws.Range("I40").Value = "LUCAS MOREIRA"
$ws.Range("G41").Value = "EDUARDO YOSHIO NAKANO"
$ws.Range("I46").Value = "EDUARDO MONTEIRO DE CASTRO GOMES"
$ws.Range("G47").Value = "DEMERSON ANDRE POLLI"

# --- Auto-fit all used columns to their (new, longer) content ---
$ws.Cells.EntireColumn.AutoFit()

# --- Restore scroll position / active selection as left by the editor ---
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D49").Select()
